# Add new columns I ("I0") and J ("IF") to the active worksheet,
# with header cells styled like the other header cells (same style as H1)
# and fill in the per-row numeric values for rows 2..76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy formatting from the existing H1 header cell so the
# new header cells pick up the same (bold/centered/bordered) style used by
# the rest of the header row, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..76 for columns I (I0) and J (IF)
$iVals = @(7,7,7,9,8,7,8,7,7,8,7,7,7,6,6,8,7,8,3,5,8,7,5,7,8,7,5,6,7,7,7,6,7,9,6,7,8,7,6,7,8,9,8,7,6,6,9,8,10,4,5,7,6,8,8,5,6,7,7,6,7,7,6,8,7,6,7,8,7,8,5,5,7,6,3)
$jVals = @(8,7,7,9,8,7,8,7,7,8,7,7,7,6,6,9,7,8,4,6,8,7,5,7,8,7,6,6,8,8,7,6,7,9,6,7,8,7,6,7,8,9,8,8,6,6,9,8,10,4,5,7,6,8,8,6,6,7,7,7,7,8,7,8,7,6,7,8,7,8,5,5,7,6,3)

for ($r = 2; $r -le 76; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
